{"js": "// docs(docs): updating metrics in briefing file\n// Adds a new \"Dia 13/09: 35min (1 dia)\" paragraph at the very end of the\n// document body, matching the formatting (Arial 12pt, 360 auto line\n// spacing, justified) already used by the other \"Dia dd/mm: ...\" entries.\n\nconst body = context.document.body;\n\n// Inserting at the \"End\" location makes Word clone the paragraph/run\n// formatting of the (current) last paragraph onto the freshly created one,\n// exactly like typing a new line at the end of the document would.\nconst newParagraph = body.insertParagraph(\"Dia 13/09: 35min (1 dia)\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# docs(docs): updating metrics in briefing file\n# Adds a new \"Dia 13/09: 35min (1 dia)\" paragraph at the very end of the\n# document body, matching the formatting (Arial 12pt, 360 auto line\n# spacing, justified) already used by the other \"Dia dd/mm: ...\" entries.\n\n$d = $word.ActiveDocument\n\n# The last paragraph in the body today is \"Dia 12/09: 3hr (1 dia)\".\n# InsertParagraphAfter() on its Range creates a brand-new paragraph right\n# after it, inheriting that paragraph's formatting (font, spacing,\n# justification) exactly the way typing Enter at the end of the document\n# would in Word.\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastParagraph.Range.InsertParagraphAfter()\n\n# The newly created paragraph is now the last one in the document; give it\n# its text.\n$newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newParagraph.Range.Text = \"Dia 13/09: 35min (1 dia)\"\n"}
